$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a numeric-looking string to a cell while keeping it as
# plain text (no leading apostrophe, no lingering number format), matching
# the original inline-string / shared-string cell content exactly.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "61.659.47"
$ws.Range("E2").Value = "  -3.87%  "
$ws.Range("D3").Value = "2.979.41"
$ws.Range("E3").Value = "  -5.07%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue $ws.Range("D5") "543.65"
$ws.Range("E5").Value = "  -4.93%  "
Set-TextValue $ws.Range("D6") "152.89"
$ws.Range("E6").Value = "  -5.68%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("D9").Value = "2.990.71"
$ws.Range("E9").Value = "  -5.04%  "
Set-TextValue $ws.Range("D10") "0.114"
$ws.Range("E10").Value = "  -1.75%  "
Set-TextValue $ws.Range("D11") "6.15"
$ws.Range("E11").Value = "  -6.59%  "
Set-TextValue $ws.Range("D12") "0.371"
$ws.Range("E12").Value = "  -3.46%  "
$ws.Range("D13").Value = "3.499.79"
$ws.Range("E13").Value = "  -5.11%  "
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("D15").Value = "61.732.86"
$ws.Range("E15").Value = "  -3.87%  "
Set-TextValue $ws.Range("D16") "23.80"
$ws.Range("E16").Value = "  -4.84%  "
$ws.Range("D17").Value = "2.981.53"
$ws.Range("E17").Value = "  -5.22%  "
$ws.Range("E18").Value = "  -4.00%  "
$ws.Range("E19").Value = "  -0.93%  "
Set-TextValue $ws.Range("D20") "12.06"
$ws.Range("E20").Value = "  -4.04%  "
Set-TextValue $ws.Range("D21") "382.33"
$ws.Range("E21").Value = "  -4.88%  "
Set-TextValue $ws.Range("D22") "6.71"
$ws.Range("E22").Value = "  -5.92%  "
$ws.Range("E23").Value = "  -0.13%  "
Set-TextValue $ws.Range("D24") "65.77"
$ws.Range("E24").Value = "  -2.92%  "
Set-TextValue $ws.Range("D25") "0.472"
$ws.Range("E25").Value = "  -2.34%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "3.102.53"
$ws.Range("E26").Value = "  -5.20%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D27") "0.190"
$ws.Range("E27").Value = "  -1.38%  "
$ws.Range("D28").Value = "0.0₃0950"
$ws.Range("E28").Value = "  -5.21%  "
$ws.Range("E29").Value = "  +0.99%  "
Set-TextValue $ws.Range("D30") "8.29"
$ws.Range("E30").Value = "  -5.76%  "
$ws.Range("E32").Value = "  -4.62%  "
Set-TextValue $ws.Range("D33") "20.56"
$ws.Range("E33").Value = "  -2.77%  "
Set-TextValue $ws.Range("D34") "161.09"
$ws.Range("E34").Value = "  +1.15%  "
Set-TextValue $ws.Range("D35") "4.71"
$ws.Range("E35").Value = "  -2.13%  "
$ws.Range("E36").Value = "  -4.48%  "
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("E38").Value = "  -4.74%  "
$ws.Range("E39").Value = "  -6.03%  "
$ws.Range("E40").Value = "  -2.95%  "
$ws.Range("D41").Value = "2.423.68"
$ws.Range("E41").Value = "  -9.08%  "
$ws.Range("E42").Value = "  -2.22%  "
Set-TextValue $ws.Range("D43") "22.23"
$ws.Range("E43").Value = "  -6.19%  "
Set-TextValue $ws.Range("D44") "0.672"
$ws.Range("E44").Value = "  -2.59%  "
Set-TextValue $ws.Range("D45") "0.0597"
$ws.Range("E45").Value = "  -2.73%  "
Set-TextValue $ws.Range("D46") "5.19"
$ws.Range("E46").Value = "  -4.61%  "
$ws.Range("E47").Value = "  -2.45%  "
$ws.Range("E48").Value = "  +0.11%  "
Set-TextValue $ws.Range("D49") "270.85"
$ws.Range("E49").Value = "  -6.24%  "
Set-TextValue $ws.Range("D50") "19.92"
$ws.Range("E50").Value = "  -5.59%  "
Set-TextValue $ws.Range("D51") "0.0957"
$ws.Range("E51").Value = "  -2.04%  "
